$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rearrange / extend the "Request" table header row (row 12) and add the
#    new "request_status" column (N12), matching the updated column order:
#    Id | type | description | client_id | user_id | sla_id | sla_expired_on |
#    sla_status | created_at | closed_at | request_status
# ---------------------------------------------------------------------------

# Copy the existing "foreign key" style (currently on I12) onto the new
# foreign-key header cells G12:I12.
$ws.Range("I12").Copy() | Out-Null
$ws.Range("G12:I12").PasteSpecial(-4122) | Out-Null

# Copy the existing "regular" header style (currently on L12) onto J12:N12
# (this also formats the brand-new N12 cell).
$ws.Range("L12").Copy() | Out-Null
$ws.Range("J12:N12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Now write the reordered header text.
$ws.Range("G12").Value = "client_id"
$ws.Range("H12").Value = "user_id"
$ws.Range("I12").Value = "sla_id"
$ws.Range("J12").Value = "sla_expired_on"
$ws.Range("K12").Value = "sla_status"
$ws.Range("L12").Value = "created_at"
$ws.Range("M12").Value = "closed_at"
$ws.Range("N12").Value = "request_status"

# ---------------------------------------------------------------------------
# 2. Extend the "Request" title merge (row 11) from D11:M11 to D11:N11.
# ---------------------------------------------------------------------------
$ws.Range("D11:M11").UnMerge() | Out-Null
$ws.Range("D11:N11").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 3. Widen columns D:N (was D:M) to the new custom width.
# ---------------------------------------------------------------------------
$ws.Range("D1:N1").EntireColumn.ColumnWidth = 13

# ---------------------------------------------------------------------------
# 4. Update the sheet view: zoom to 190% and move the active selection to H13.
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Zoom = 190
$ws.Range("H13").Select() | Out-Null
